{"js": "const body = context.document.body;\n\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n    await context.sync();\n  }\n}\n\n// 1. Title heading\nawait replaceOnce(\n  \"Play Cash Falls Island Bounty for Free - Review\",\n  \"Play Cash Falls Island Bounty for Free\"\n);\n\n// 2-5. \"What we like\" bullet list\nawait replaceOnce(\"High number of ways to win\", \"Beautiful Caribbean beach theme\");\nawait replaceOnce(\"Jackpot opportunities\", \"Large number of ways to win\");\nawait replaceOnce(\"Multiple bonus features\", \"Exciting gameplay features like Cash Falls and Jackpots\");\nawait replaceOnce(\"Turbo and Autoplay options\", \"Opportunity for high payouts\");\n\n// 6. \"What we don't like\" bullet list\nawait replaceOnce(\"RTP is lower than some other slot games\", \"Not available for free play\");\n\n// 7. Bold \"Play ... Review\" line near the end\nawait replaceOnce(\n  \"Play Cash Falls Island Bounty for Free - Review\",\n  \"Play Cash Falls Island Bounty for Free\"\n);\n\n// 8. Italic summary line\nawait replaceOnce(\n  \"Read our neutral review of Cash Falls Island Bounty. Play this online slot game for free and find out about its gameplay features, symbols, and much more.\",\n  \"Read our review of Cash Falls Island Bounty and play this exciting slot game for free.\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-OneText($FindText, $ReplaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.Text = $ReplaceText\n    $find.Execute(\n        $FindText,        # FindText\n        $true,             # MatchCase\n        $false,            # MatchWholeWord\n        $false,            # MatchWildcards\n        $false,            # MatchSoundsLike\n        $false,            # MatchAllWordForms\n        $true,             # Forward\n        1,                 # Wrap: wdFindContinue\n        $false,            # Format\n        $ReplaceText,      # ReplaceWith\n        1                  # Replace: wdReplaceOne\n    ) | Out-Null\n}\n\n# 1. Title heading\nReplace-OneText \"Play Cash Falls Island Bounty for Free - Review\" \"Play Cash Falls Island Bounty for Free\"\n\n# 2-5. \"What we like\" bullet list\nReplace-OneText \"High number of ways to win\" \"Beautiful Caribbean beach theme\"\nReplace-OneText \"Jackpot opportunities\" \"Large number of ways to win\"\nReplace-OneText \"Multiple bonus features\" \"Exciting gameplay features like Cash Falls and Jackpots\"\nReplace-OneText \"Turbo and Autoplay options\" \"Opportunity for high payouts\"\n\n# 6. \"What we don't like\" bullet list\nReplace-OneText \"RTP is lower than some other slot games\" \"Not available for free play\"\n\n# 7. Bold \"Play ... Review\" line near the end\nReplace-OneText \"Play Cash Falls Island Bounty for Free - Review\" \"Play Cash Falls Island Bounty for Free\"\n\n# 8. Italic summary line\nReplace-OneText \"Read our neutral review of Cash Falls Island Bounty. Play this online slot game for free and find out about its gameplay features, symbols, and much more.\" \"Read our review of Cash Falls Island Bounty and play this exciting slot game for free.\"\n"}
